$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose column-B cell must be cleared entirely (no cell left in the row)
$clearRows = @(4, 5, 6, 11, 12, 13, 14, 15)
foreach ($r in $clearRows) {
    $ws.Range("B$r").ClearContents()
}

# Rows B7, B8, B9: "entity id, UNC, UC" -> "entity id, "
# Row B10: "sensor number, entity id, UNC, UC" -> "entity id, " (also drops "sensor number, ")
foreach ($r in @(7, 8, 9, 10)) {
    $ws.Range("B$r").Value = "entity id, "
}

# Rows B17..B47: "sensor number, UNC, UC" -> "sensor number, "
for ($r = 17; $r -le 47; $r++) {
    $ws.Range("B$r").Value = "sensor number, "
}

# Rows B48, B49: "Sensor name, sensor number, UNC, UC" -> "Sensor name, sensor number, "
foreach ($r in @(48, 49)) {
    $ws.Range("B$r").Value = "Sensor name, sensor number, "
}

# Rows B50..B62: "sensor number, UNC, UC" -> "sensor number, "
for ($r = 50; $r -le 62; $r++) {
    $ws.Range("B$r").Value = "sensor number, "
}
